# Append-run update for the "ランサーズ" (lancers) sheet:
# - Scrape timestamp moves from 2026-01-14 02:04:53 to 2026-01-14 06:30:40
# - The result set shrinks from 17 listings (rows 2-18) to 6 listings (rows 2-7)
# - A few columns get narrower
# - Hyperlinks in column F need to point at the new listing URLs

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Drop the rows that no longer exist in the new result set -----------
# (rows 8-18 of the old sheet). This also shifts the dimension down to A1:H7
# automatically.
$ws.Range("A8:H18").EntireRow.Delete()

# --- 2. Column width tweaks --------------------------------------------------
# Excel's ColumnWidth property pads by 5/6 of a character vs. the stored
# OOXML <col width>, so subtract that back out to land on the exact target
# widths (50 / 28 / 12).
$ws.Columns.Item(2).ColumnWidth = 50 - 5/6
$ws.Columns.Item(4).ColumnWidth = 28 - 5/6
$ws.Columns.Item(8).ColumnWidth = 12 - 5/6

# --- 3. Rebuild the 6 data rows with the new scrape's content ---------------
$rows = @(
    @{ Row=2;  A="2026-01-14 06:30:40"; B="AI企画書作成システムの「見積書作成」をご支援いただける制作会社/エンジニア募集(発注確約なし)"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";  E="期限情報なし"; F="https://www.lancers.jp/work/detail/5470737"; G=313; H="🔥AI,Ai" },
    @{ Row=3;  A="2026-01-14 06:30:40"; B="進行管理およびチームディレクションを担当"; C="システム開発"; D="~ 5,000 円 / 固定";              E="期限情報なし"; F="https://www.lancers.jp/work/detail/5418064"; G=30;  H="◇管理" },
    @{ Row=4;  A="2026-01-14 06:30:40"; B="Rubyの暗号化機能のPHP化"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";  E="期限情報なし"; F="https://www.lancers.jp/work/detail/5470623"; G=28;  H="○PHP" },
    @{ Row=5;  A="2026-01-14 06:30:40"; B="金融機関の入出金伝票印刷システム構築依頼"; C="システム開発"; D="20,000 円 ~ 50,000 円 / 固定";  E="期限情報なし"; F="https://www.lancers.jp/work/detail/5470403"; G=28;  H=$null },
    @{ Row=6;  A="2026-01-14 06:30:40"; B="《長期レギュラー》公的機関Web運用の要となる、ディレクター募集"; C="システム開発"; D="200,000 円 ~ 300,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5470150"; G=18;  H=$null },
    @{ Row=7;  A="2026-01-14 06:30:40"; B="【フリーランス必見】エンジニア支援サービスのご紹介!"; C="システム開発"; D="10,000 円 ~ 20,000 円 / 固定"; E="期限情報なし"; F="https://www.lancers.jp/work/detail/5470726"; G=10;  H=$null }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Cells.Item($n, 1).Value2 = $r.A
    $ws.Cells.Item($n, 2).Value2 = $r.B
    $ws.Cells.Item($n, 3).Value2 = $r.C
    $ws.Cells.Item($n, 4).Value2 = $r.D
    $ws.Cells.Item($n, 5).Value2 = $r.E
    $ws.Cells.Item($n, 6).Value2 = $r.F
    $ws.Cells.Item($n, 7).Value2 = $r.G
    if ($r.H) {
        $ws.Cells.Item($n, 8).Value2 = $r.H
    } else {
        $ws.Cells.Item($n, 8).ClearContents()
    }
}

# --- 4. Rebuild hyperlinks for column F (rows 2-7) --------------------------
# Wipe the stale hyperlink collection (any cell-scoped Hyperlinks.Delete()
# clears the whole sheet collection in this runtime) then add exactly the six
# links the new result set needs, pointing at the refreshed URLs.
$ws.Range("F2").Hyperlinks.Delete()

foreach ($r in $rows) {
    $cell = $ws.Cells.Item($r.Row, 6)
    $ws.Hyperlinks.Add($cell, $r.F) | Out-Null
    $cell.Style = "Hyperlink"
}
